$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("November")

$ws.Range("B31").Value = 1
$ws.Range("C31").Value = "UART problem nachgehen"
$ws.Range("D31").Value = "19:00 - 20:00"

$ws.Range("B32").Value = 2
$ws.Range("C32").Value = "UART problem nachgehen"
$ws.Range("D32").Value = "17:00 - 19:00"

$ws.Range("B33").Value = 2.5
$ws.Range("D33").Value = "14:00 - 16:30"
$ws.Range("C33").Value = "Teammeating, zusammenführen UART und RF"

$ws.Range("D32").Select()
